$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants por usuario")

$rows = @(
    @{A="GRANT DELETE, INSERT, SELECT, UPDATE ON TABLE tratamento.tb_hstr_pnel_solic_trtmto TO"; D="GRANT DELETE, INSERT, SELECT, UPDATE ON TABLE tratamento.tb_hstr_pnel_solic_trtmto TO aline ;"},
    @{A="GRANT ALL ON SEQUENCE tratamento.sq_hstr_pnel_solic_trtmto TO "; D="GRANT ALL ON SEQUENCE tratamento.sq_hstr_pnel_solic_trtmto TO  aline ;"},
    @{A="GRANT DELETE, INSERT, SELECT, UPDATE ON TABLE tratamento.tb_c_pcnt TO"; D="GRANT DELETE, INSERT, SELECT, UPDATE ON TABLE tratamento.tb_c_pcnt TO aline ;"},
    @{A="GRANT SELECT ON  TABLE tratamento.vw_painel_trtmto TO"; D="GRANT SELECT ON  TABLE tratamento.vw_painel_trtmto TO aline ;"}
)

$r = 39
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = "aline"
    $ws.Cells.Item($r, 2).Style = $ws.Cells.Item($r - 1, 2).Style
    $ws.Cells.Item($r, 3).Value = ";"
    $ws.Cells.Item($r, 4).Formula = "=A$r&`" `"&B$r&`" `"&C$r"
    $r++
}

$ws.Range("A39").Select()
